# Lisää Tielupien haku roolit-exceliin
# Insert a new row into the "Oikeudet" sheet (just above the "Ilmoitukset" row,
# i.e. at row 84) for the new "Tieluvat" / "Tielupien haku" entry, duplicating
# the surrounding row's formatting and R*/W* access-pattern values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")

# Insert a fresh row at 84 (this pushes the former row 84 "Ilmoitukset" down to 85).
$ws.Rows.Item(84).Insert()

# Copy the (now shifted-down) row 85 formatting into the new blank row 84, so the
# new row matches the look of the table it is being inserted into.
$ws.Range("A85:Y85").Copy()
$ws.Range("A84:Y84").PasteSpecial(-4122)

# Fill in the values for the new "Tieluvat" row.
$ws.Cells.Item(84, 1).Value2 = "Tieluvat"
$ws.Cells.Item(84, 2).Value2 = "Tielupien haku"

$ws.Cells.Item(84, 4).Value2  = "R*,W*"
$ws.Cells.Item(84, 5).Value2  = "R*"
$ws.Cells.Item(84, 6).Value2  = "R*"
$ws.Cells.Item(84, 7).Value2  = "R*"
$ws.Cells.Item(84, 8).Value2  = "R*"
$ws.Cells.Item(84, 9).Value2  = "R*"
$ws.Cells.Item(84, 10).Value2 = "R"
$ws.Cells.Item(84, 11).Value2 = "R*,W*"
$ws.Cells.Item(84, 12).Value2 = "R*"
$ws.Cells.Item(84, 13).Value2 = "R"
$ws.Cells.Item(84, 14).Value2 = "R*"
$ws.Cells.Item(84, 15).Value2 = "R*"
$ws.Cells.Item(84, 16).Value2 = "R*"
$ws.Cells.Item(84, 17).Value2 = "R"
$ws.Cells.Item(84, 18).Value2 = "R+,W+"
$ws.Cells.Item(84, 19).Value2 = "R,W"
$ws.Cells.Item(84, 20).Value2 = "R+"
$ws.Cells.Item(84, 21).Value2 = "R,W"
$ws.Cells.Item(84, 22).Value2 = "R+"
$ws.Cells.Item(84, 23).Value2 = "R,W"
$ws.Cells.Item(84, 24).Value2 = "R,W"

# Column C stays blank, matching the rest of the table.
$ws.Cells.Item(84, 3).Value2 = ""

# Keep the filter range / selection roughly in sync with the now-larger table.
$ws.Range("D80").Select()
